# Updates cryptos list figures (price + 1h volume change) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "25.928.33"
$cell.ClearFormats()
$ws.Range("E2").Value = "  -0.50%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.638.36"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +0.05%  "

$ws.Range("E4").Value = "  +0.99%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "214.54"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -0.31%  "

$ws.Range("E6").Value = "  +0.56%  "

$ws.Range("E7").Value = "  +0.75%  "

$ws.Range("E8").Value = "  +0.73%  "

$ws.Range("E9").Value = "  -1.25%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "19.60"
$cell.ClearFormats()
$ws.Range("E10").Value = "  -0.77%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0794"
$cell.ClearFormats()

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.864.86"
$cell.ClearFormats()
$ws.Range("E12").Value = "  -0.19%  "

$ws.Range("E13").Value = "  -0.19%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "1.607.89"
$cell.ClearFormats()
$ws.Range("E14").Value = "  -1.90%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.543"
$cell.ClearFormats()
$ws.Range("E15").Value = "  -1.71%  "

$ws.Range("E16").Value = "  -0.62%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "62.52"
$cell.ClearFormats()
$ws.Range("E17").Value = "  -1.31%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "25.957.39"
$cell.ClearFormats()
$ws.Range("E18").Value = "  -0.52%  "

$ws.Range("E19").Value = "  +0.97%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "193.43"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +0.25%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "4.38"
$cell.ClearFormats()
$ws.Range("E21").Value = "  -1.58%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "9.91"
$cell.ClearFormats()
$ws.Range("E22").Value = "  -0.96%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "6.27"
$cell.ClearFormats()
$ws.Range("E23").Value = "  -1.99%  "

$ws.Range("E24").Value = "  +0.67%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "143.99"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +1.61%  "

$ws.Range("E26").Value = "  +0.96%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.126"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +2.74%  "

$ws.Range("E28").Value = "  -0.64%  "

$ws.Range("E29").Value = "  -0.89%  "

$ws.Range("E30").Value = "  +0.08%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.0499"
$cell.ClearFormats()
$ws.Range("E31").Value = "  +1.09%  "

$ws.Range("E32").Value = "  -0.99%  "

$ws.Range("E33").Value = "  -0.56%  "

$ws.Range("E34").Value = "  -3.48%  "

$ws.Range("E35").Value = "  +2.15%  "

$ws.Range("E36").Value = "  -0.63%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.138.24"
$cell.ClearFormats()
$ws.Range("E37").Value = "  -1.12%  "

$ws.Range("E38").Value = "  -0.07%  "

$ws.Range("E39").Value = "  -1.28%  "

$ws.Range("E40").Value = "  +0.14%  "

$ws.Range("E41").Value = "  +0.82%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "99.38"
$cell.ClearFormats()
$ws.Range("E42").Value = "  -1.00%  "

$ws.Range("E43").Value = "  -0.14%  "

$ws.Range("E44").Value = "  -3.90%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.774.17"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -0.36%  "

$ws.Range("E46").Value = "  +5.23%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "56.55"
$cell.ClearFormats()
$ws.Range("E47").Value = "  +1.71%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.0529"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +2.13%  "

$ws.Range("E49").Value = "  -1.08%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "7.67"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +0.66%  "

$ws.Range("E51").Value = "  +0.14%  "
